$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "max" column (old column C). Shift "prediction" (old D) into C,
# and "rejection-f" (old E) into D, then delete the now-unused old E column.
$ws.Columns.Item(3).Delete()

# Header row stays the same text, just confirm values after shift.
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"

# Update data row values per new outputs.
$ws.Range("B2").Value = -11210.93793505753
$ws.Range("C2").Value = "g__RUG033"
$ws.Range("D2").Value = "g__RUG033(reject)"
